$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 126, pushing the existing rows 126-129 down to 128-131.
$ws.Rows.Item(126).Resize(2).Insert()

# Row 126: new weekly record (date 2021-09-09 = serial 44448), Zafiro rojo
$ws.Range("A126").Value = 7
$ws.Range("B126").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C126").Value = "Ñuble"
$ws.Range("D126").Value = 44448
$ws.Range("E126").Value = 16
$ws.Range("F126").Value = 100112002
$ws.Range("G126").Value = "Pimiento"
$ws.Range("H126").Value = "Zafiro rojo"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 160
$ws.Range("K126").Value = 44000
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = 44500
$ws.Range("N126").Value = "$/caja 15 kilos"
$ws.Range("O126").Value = "Región de Arica y Parinacota"
$ws.Range("P126").Value = 2967
$ws.Range("Q126").Value = 15
$ws.Range("R126").Value = "Hortaliza"

# Row 127: new weekly record (date 2021-09-09 = serial 44448), Zafiro verde
$ws.Range("A127").Value = 7
$ws.Range("B127").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C127").Value = "Ñuble"
$ws.Range("D127").Value = 44448
$ws.Range("E127").Value = 16
$ws.Range("F127").Value = 100112002
$ws.Range("G127").Value = "Pimiento"
$ws.Range("H127").Value = "Zafiro verde"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 160
$ws.Range("K127").Value = 41000
$ws.Range("L127").Value = 42000
$ws.Range("M127").Value = 41500
$ws.Range("N127").Value = "$/caja 15 kilos"
$ws.Range("O127").Value = "Región de Arica y Parinacota"
$ws.Range("P127").Value = 2767
$ws.Range("Q127").Value = 15
$ws.Range("R127").Value = "Hortaliza"
